$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P2").Value = 1.24
$ws.Range("F3").Value = 2.86
$ws.Range("G3").Value = 3.15
$ws.Range("H3").Value = 2.48
$ws.Range("I3").Value = 2.72
$ws.Range("J3").Value = 3.35
$ws.Range("K3").Value = 3.75
$ws.Range("P3").Value = 1.88
$ws.Range("Q3").Value = 1.97
$ws.Range("F4").Value = 3.45
$ws.Range("G4").Value = 3.9
$ws.Range("H4").Value = 2.42
$ws.Range("I4").Value = 2.7
$ws.Range("K4").Value = 3.1
$ws.Range("N4").Value = 2.48
$ws.Range("O4").Value = 1.55
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.7
$ws.Range("T4").Value = 2.06
$ws.Range("U4").Value = 1.74
$ws.Range("X4").Value = 9.6
$ws.Range("Y4").Value = 9
$ws.Range("Z4").Value = 18
$ws.Range("AA4").Value = 50
$ws.Range("AB4").Value = 11
$ws.Range("AC4").Value = 7.2
$ws.Range("AD4").Value = 13.5
$ws.Range("AE4").Value = 46
$ws.Range("AF4").Value = 28
$ws.Range("AG4").Value = 20
$ws.Range("AH4").Value = 28
$ws.Range("AI4").Value = 70
$ws.Range("AJ4").Value = 90
$ws.Range("AK4").Value = 65
$ws.Range("AL4").Value = 110
$ws.Range("AM4").Value = 230
$ws.Range("AN4").Value = 100
$ws.Range("AO4").Value = 55
$ws.Range("F5").Value = 3.2
$ws.Range("G5").Value = 3.6
$ws.Range("H5").Value = 2.52
$ws.Range("I5").Value = 2.86
$ws.Range("J5").Value = 2.88
$ws.Range("K5").Value = 3.35
$ws.Range("P5").Value = 1.54
$ws.Range("Q5").Value = 2.48
$ws.Range("F6").Value = 3.15
$ws.Range("G6").Value = 3.7
$ws.Range("H6").Value = 2.3
$ws.Range("I6").Value = 2.6
$ws.Range("J6").Value = 3.1
$ws.Range("K6").Value = 3.65
$ws.Range("P6").Value = 1.73
$ws.Range("Q6").Value = 2.1
$ws.Range("F8").Value = 5.7
$ws.Range("G8").Value = 6.8
$ws.Range("H8").Value = 1.7
$ws.Range("I8").Value = 1.83
$ws.Range("J8").Value = 3.45
$ws.Range("K8").Value = 3.95
$ws.Range("P8").Value = 1.66
$ws.Range("Q8").Value = 2.26
$ws.Range("F9").Value = 3.7
$ws.Range("G9").Value = 4.2
$ws.Range("I9").Value = 2.2
$ws.Range("J9").Value = 3.5
$ws.Range("P9").Value = 1.92
$ws.Range("F10").Value = 3.5
$ws.Range("G10").Value = 3.7
$ws.Range("H10").Value = 2.22
$ws.Range("J10").Value = 3.55
$ws.Range("K10").Value = 3.65
$ws.Range("M10").Value = 1.07
$ws.Range("P10").Value = 1.84
$ws.Range("Q10").Value = 2
$ws.Range("T10").Value = 1.78
$ws.Range("U10").Value = 2.08
$ws.Range("X10").Value = 14
$ws.Range("Y10").Value = 10
$ws.Range("Z10").Value = 17.5
$ws.Range("AA10").Value = 32
$ws.Range("AB10").Value = 14
$ws.Range("AC10").Value = 8.4
$ws.Range("AD10").Value = 13.5
$ws.Range("AE10").Value = 32
$ws.Range("AF10").Value = 27
$ws.Range("AG10").Value = 16
$ws.Range("AH10").Value = 19
$ws.Range("AL10").Value = 60
$ws.Range("AM10").Value = 130
$ws.Range("AO10").Value = 20
$ws.Range("F11").Value = 4.5
$ws.Range("G11").Value = 5.2
$ws.Range("H11").Value = 1.89
$ws.Range("I11").Value = 1.97
$ws.Range("J11").Value = 3.55
$ws.Range("P11").Value = 1.89
$ws.Range("Q11").Value = 1.98
$ws.Range("T11").Value = 1.84
$ws.Range("U11").Value = 2.04
$ws.Range("AA11").Value = 980
$ws.Range("AE11").Value = 980
$ws.Range("AF11").Value = 980
$ws.Range("AG11").Value = 23
$ws.Range("AI11").Value = 980
$ws.Range("AJ11").Value = 120
$ws.Range("AK11").Value = 65
$ws.Range("AL11").Value = 70
$ws.Range("AM11").Value = 120
$ws.Range("AN11").Value = 75
$ws.Range("F12").Value = 1.31
$ws.Range("H12").Value = 11.5
$ws.Range("I12").Value = 13.5
$ws.Range("J12").Value = 5.7
$ws.Range("K12").Value = 6.2
$ws.Range("N12").Value = 5.2
$ws.Range("P12").Value = 2.44
$ws.Range("Q12").Value = 1.6
$ws.Range("R12").Value = 1.58
$ws.Range("S12").Value = 2.5
$ws.Range("T12").Value = 2.06
$ws.Range("U12").Value = 1.83
$ws.Range("X12").Value = 26
$ws.Range("Z12").Value = 130
$ws.Range("AA12").Value = 520
$ws.Range("AB12").Value = 9.6
$ws.Range("AC12").Value = 13.5
$ws.Range("AE12").Value = 220
$ws.Range("AF12").Value = 8.6
$ws.Range("AG12").Value = 11
$ws.Range("AI12").Value = 170
$ws.Range("AJ12").Value = 10.5
$ws.Range("AM12").Value = 180
$ws.Range("AN12").Value = 4.9
$ws.Range("H13").Value = 5.6
$ws.Range("I13").Value = 6.4
$ws.Range("J13").Value = 3.9
$ws.Range("K13").Value = 4.3
$ws.Range("P13").Value = 2.04
$ws.Range("Q13").Value = 1.84
$ws.Range("T13").Value = 1.83
$ws.Range("U13").Value = 2.04
$ws.Range("X13").Value = 17.5
$ws.Range("AA13").Value = 170
$ws.Range("AB13").Value = 9
$ws.Range("AC13").Value = 9.199999999999999
$ws.Range("AF13").Value = 11
$ws.Range("AM13").Value = 130
$ws.Range("AN13").Value = 9.800000000000001
$ws.Range("F15").Value = 6
$ws.Range("G15").Value = 6.8
$ws.Range("J15").Value = 4.5
$ws.Range("K15").Value = 5
$ws.Range("N15").Value = 5
$ws.Range("P15").Value = 2.38
$ws.Range("Q15").Value = 1.64
$ws.Range("R15").Value = 1.56
$ws.Range("S15").Value = 2.56
$ws.Range("T15").Value = 1.77
$ws.Range("U15").Value = 2.16
$ws.Range("X15").Value = 25
$ws.Range("Y15").Value = 10.5
$ws.Range("Z15").Value = 11
$ws.Range("AA15").Value = 15.5
$ws.Range("AC15").Value = 11
$ws.Range("AD15").Value = 10.5
$ws.Range("AG15").Value = 27
$ws.Range("AJ15").Value = 190
$ws.Range("AO15").Value = 6.8
$ws.Range("F16").Value = 2.24
$ws.Range("G16").Value = 2.38
$ws.Range("I16").Value = 3.95
$ws.Range("J16").Value = 3.3
$ws.Range("K16").Value = 3.4
$ws.Range("P16").Value = 1.78
$ws.Range("G17").Value = 1.48
$ws.Range("H17").Value = 8.199999999999999
$ws.Range("I17").Value = 10.5
$ws.Range("K17").Value = 5.1
$ws.Range("P17").Value = 2.06
$ws.Range("Q17").Value = 1.87
$ws.Range("F19").Value = 2.96
$ws.Range("G19").Value = 3.25
$ws.Range("H19").Value = 2.62
$ws.Range("I19").Value = 2.8
$ws.Range("J19").Value = 3.25
$ws.Range("K19").Value = 3.3
$ws.Range("P19").Value = 1.68
$ws.Range("Q19").Value = 2.3
$ws.Range("F20").Value = 2.16
$ws.Range("I20").Value = 4.1
$ws.Range("K20").Value = 3.5
$ws.Range("P20").Value = 1.7
$ws.Range("Q20").Value = 2.26
